# Aggiunta tabelle crv in file snapshot
# Adds 6 new rows (50-55) to the "Snapshot" sheet, describing the "crv" / "curato"
# snapshot tables tcr01_forniture_gas, tcr01_forniture_pwr, tcr01_mapping_tcr_pwr,
# tcr01_ricavi_gas, tcr01_ricavi_pwr, tcr01_mapping_tcr_gas - mirroring the existing
# rows 28-43 layout/format (same A/B/D/E/F pattern, only C changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snapshot")

# 1) Clone formatting (style + border + row look) of the existing "tcr01" block
#    (row 28) down onto the new rows in one shot - Excel tiles a single-row
#    source across a multi-row destination.
$ws.Range("A28:I28").Copy($ws.Range("A50:I55"))

# 2) Fill in the actual content for each new row.
$newRows = @(
    @{ Row = 50; C = "tcr01_forniture_gas" },
    @{ Row = 51; C = "tcr01_forniture_pwr" },
    @{ Row = 52; C = "tcr01_mapping_tcr_pwr" },
    @{ Row = 53; C = "tcr01_ricavi_gas" },
    @{ Row = 54; C = "tcr01_ricavi_pwr" },
    @{ Row = 55; C = "tcr01_mapping_tcr_gas" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = "crv"
    $ws.Range("B$r").Value = "curato"
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = "snapshot"
    $ws.Range("E$r").Value = "dat_startDate"
    $ws.Range("F$r").Value = "timestamp"
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = ""
    $ws.Rows.Item($r).RowHeight = 15
}

# 3) Restore the view state: scrolled so row 25 is at the top, with F53 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("F53").Select()
